$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comparisons")

$ws.Range("C2").Value = 15.796417236328125
$ws.Range("E2").Value = 2229.0

$ws.Range("C4").Value = 4.249732494354248
$ws.Range("E4").Value = 437.0

$ws.Range("C5").Value = 2.2626311779022217
$ws.Range("E5").Value = 367.0

$ws.Range("C6").Value = 2.4519267082214355
$ws.Range("E6").Value = 370.0

$ws.Range("C7").Value = 2.2626311779022217
$ws.Range("E7").Value = 367.0

$ws.Range("C8").Value = 2.3265609741210938
$ws.Range("E8").Value = 355.0

$ws.Range("C9").Value = 2.242934226989746
$ws.Range("E9").Value = 333.0
